$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C = y_0_forecast, Column E = y_1_forecast
$ws.Range("C2").Value = 2.071001150990881
$ws.Range("E2").Value = 5.125398993052044

$ws.Range("C3").Value = 0.8787496612562951
$ws.Range("E3").Value = 1.013823151053095

$ws.Range("C4").Value = 2.723861837742825
$ws.Range("E4").Value = 4.356912452939476

$ws.Range("C5").Value = 2.431458940167008
$ws.Range("E5").Value = 5.259925231829898

$ws.Range("C6").Value = 1.447930496829541
$ws.Range("E6").Value = 0.7814198158142105

$ws.Range("C7").Value = 0.3494148569448852
$ws.Range("E7").Value = -0.5632089058212553

$ws.Range("C8").Value = 0.6742451383205061
$ws.Range("E8").Value = 1.713290556413605

$ws.Range("C9").Value = 2.178094576990031
$ws.Range("E9").Value = 2.743122680804988

$ws.Range("C10").Value = 2.623600596229347
$ws.Range("E10").Value = 1.891220645890002

$ws.Range("C11").Value = 1.715791310593229
$ws.Range("E11").Value = 1.687339605296523

$ws.Range("C12").Value = 0.3317798769387315
$ws.Range("E12").Value = -2.079848588862154

$ws.Range("C13").Value = 0.8934982674867697
$ws.Range("E13").Value = -1.194610791899997

$ws.Range("C14").Value = 1.979074033580819
$ws.Range("E14").Value = 0.8024032015999882

$ws.Range("C15").Value = 2.552476296061434
$ws.Range("E15").Value = 3.086122033237126

$ws.Range("C16").Value = 0.3901728183783204
$ws.Range("E16").Value = 1.906002353653125

$ws.Range("C17").Value = -1.906744368254853
$ws.Range("E17").Value = 9.556389850000446

$ws.Range("C18").Value = 1.152760694685062
$ws.Range("E18").Value = 0.6444718444275521

$ws.Range("C19").Value = 2.287212358310953
$ws.Range("E19").Value = 1.801540135156521
